$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting from row 11 down to row 12 so the new row matches
# the existing styling (number formats per column).
$ws.Range("A11:T11").Copy() | Out-Null
$ws.Range("A12:T12").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A12").Value = "v10-3200"
$ws.Range("B12").Value = 60
$ws.Range("C12").Value = "Fixed severe bugs in state mapping"
$ws.Range("D12").Value = 20
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 1620
$ws.Range("J12").Value = 1620
$ws.Range("K12").Value = 1620
$ws.Range("L12").Value = 154
$ws.Range("M12").Value = 154
$ws.Range("N12").Value = 154
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 0
$ws.Range("R12").Value = 97.3
$ws.Range("S12").Value = 97.3
$ws.Range("T12").Value = 97.3

# Update the selected cell, matching the diff's recorded view state.
$ws.Range("I20").Select() | Out-Null
